$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.475.75"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.801.21"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.51"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.41"
$ws.Range("E8").Value = "  +18.00%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0665"
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0995"
$ws.Range("E11").Value = "  +3.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.062.88"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.802.61"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.91"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.451.89"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.40"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.09"
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.14"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.36"
$ws.Range("E23").Value = "  +6.83%  "
$ws.Range("E24").Value = "  -2.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.15"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.64"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.38"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.79"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "87.35"
$ws.Range("E35").Value = "  +7.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.646"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.317.96"
$ws.Range("E37").Value = "  -3.37%  "
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0188"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.75"
$ws.Range("E40").Value = "  +11.58%  "
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  +4.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.43"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.935"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0518"
$ws.Range("E46").Value = "  +3.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.963.52"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.79"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "100.40"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("E51").Value = "  +0.53%  "
